# Features List.xlsx update
# - Rewrites the "Gamify this thing" details cell (D5) with the expanded,
#   partially-highlighted announcer-line notes (adds a "pause menu" bullet and
#   extends the closing "any thoughts about art?" line).
# - Removes the now-obsolete "Visual Upgrades" row (row 13) entirely.
# - Grows row 5 to fit the longer text and resets the saved scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Clear out the old "Visual Upgrades" row -------------------------------
$ws.Range("B13").ClearContents()
$ws.Range("C13").ClearContents()
$ws.Range("D13").ClearContents()

# --- Rewrite the "Gamify this thing" details cell --------------------------
$cell = $ws.Range("D5")

$newText = "This thing is a prototype currently, lets move to make it more gamified`n" + `
"       add a menu system with buttons to start a game, an option to quit to menu from game, an option to quit the game entirely, etc. Add whatever makes sense here`n" + `
"       get us to a place where we can generate builds and they'll actually be playable`n" + `
"       a set up to run multiple scenes so that I can customize the ai behaviour and pawn skill sliders for each differen scene`n" + `
"       potentially an options menu, allow players to set things like sound volume`n" + `
"       potentially a credits page that has links to both our portfolio websites`n" + `
"       pause menu that has the option to include game statistics (I am tracking these already) and to quit the game`n" + `
"       any thoughts about art? making this thing look better? any other stuff we need to put this thing out there?"

$cell.Value2 = $newText

# Highlight the same three "bullet" segments in green (RGB 146,208,80 -> VBA
# BGR long 5296274), matching the existing formatting convention used
# throughout this sheet for in-progress/partially-implemented notes.
$greenColor = 5296274
$cell.Characters(80, 156).Font.Color = $greenColor
$cell.Characters(244, 80).Font.Color = $greenColor
$cell.Characters(331, 119).Font.Color = $greenColor

# --- Row height grew to fit the extra lines of text -------------------------
$ws.Rows.Item(5).RowHeight = 120

# --- Reset the saved scroll position (no more topLeftCell override) --------
$excel.ActiveWindow.ScrollRow = 1
